$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.175.17"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.587.96"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.26%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0602"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.38%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.00"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.44%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0844"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "1.811.63"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "1.584.84"
$ws.Range("E13").Value = "  -0.36%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.01"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.510"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.66%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "63.46"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "26.180.94"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("E18").Value = "  -0.54%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.20%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "213.67"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("E21").Value = "  -0.11%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.23"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.71%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.12"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -1.77%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.112"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.24%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "15.03"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.41%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0493"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").Value = "1.410.68"
$ws.Range("E33").Value = "  +8.11%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.72%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.42"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("E38").Value = "  -1.73%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.821"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.88%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.87"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.77%  "
$ws.Range("E41").Value = "  -0.14%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.946"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -13.43%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.14"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.30%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.761"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").Value = "1.723.26"
$ws.Range("E45").Value = "  -0.04%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "60.96"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.28%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "85.89"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  -0.88%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0955"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.79%  "
